# Apply updated Price values (column D) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 9.630000000000001
$ws.Range("D3").Value = 13.86
$ws.Range("D4").Value = 17.06
$ws.Range("D5").Value = 9.58
$ws.Range("D6").Value = 8.73
$ws.Range("D8").Value = 0.86
$ws.Range("D9").Value = 1.43
$ws.Range("D11").Value = 3.28
$ws.Range("D14").Value = 0.3
$ws.Range("D15").Value = 0.53
$ws.Range("D16").Value = 0.13
$ws.Range("D18").Value = 0.58
$ws.Range("D19").Value = 7.97
$ws.Range("D20").Value = 2.59
$ws.Range("D21").Value = 2.99
$ws.Range("D22").Value = 0.76
$ws.Range("D24").Value = 32.43
$ws.Range("D25").Value = 6.43
$ws.Range("D26").Value = 19.22
$ws.Range("D28").Value = 1.01
$ws.Range("D29").Value = 4.73
$ws.Range("D30").Value = 3.93
$ws.Range("D32").Value = 0.77
$ws.Range("D34").Value = 0.91
$ws.Range("D35").Value = 15.28
$ws.Range("D36").Value = 4.69
$ws.Range("D37").Value = 6.79
$ws.Range("D38").Value = 7.89
$ws.Range("D39").Value = 3.86
$ws.Range("D40").Value = 0.7
$ws.Range("D42").Value = 56.96
$ws.Range("D43").Value = 0.87
$ws.Range("D44").Value = 0.87
$ws.Range("D45").Value = 4.82
$ws.Range("D46").Value = 0.82
$ws.Range("D47").Value = 2.9
$ws.Range("D48").Value = 7.52
$ws.Range("D49").Value = 1.73
$ws.Range("D51").Value = 4.93
$ws.Range("D52").Value = 1.14
$ws.Range("D55").Value = 1.25
$ws.Range("D57").Value = 1.53
$ws.Range("D58").Value = 1.53
$ws.Range("D59").Value = 1.97
$ws.Range("D60").Value = 1.49
$ws.Range("D61").Value = 2.75
$ws.Range("D63").Value = 1.69
$ws.Range("D64").Value = 4.76
$ws.Range("D65").Value = 26.92
$ws.Range("D69").Value = 3.45
$ws.Range("D70").Value = 4.09
$ws.Range("D71").Value = 0.53
$ws.Range("D72").Value = 1.01
$ws.Range("D73").Value = 2.67
$ws.Range("D74").Value = 2.87
$ws.Range("D75").Value = 6.95
$ws.Range("D79").Value = 6
$ws.Range("D80").Value = 3.86
$ws.Range("D81").Value = 4.35
$ws.Range("D83").Value = 1.38
$ws.Range("D84").Value = 11.83
$ws.Range("D85").Value = 6.7
$ws.Range("D86").Value = 17.02
$ws.Range("D87").Value = 1.19
$ws.Range("D88").Value = 10.93
$ws.Range("D89").Value = 9.470000000000001
$ws.Range("D90").Value = 3.62
$ws.Range("D91").Value = 1.81
$ws.Range("D93").Value = 5.95
$ws.Range("D97").Value = 1.93
$ws.Range("D100").Value = 1.09
